# Update "想去人数" (interest count) and "最低票价" (lowest ticket price)
# figures on the 展览 and 全部类型 sheets to the freshly scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 3769
$ws1.Range("F9").Value  = 112
$ws1.Range("F12").Value = 81
$ws1.Range("F15").Value = 895
$ws1.Range("G16").Value = 55
$ws1.Range("F17").Value = 226
$ws1.Range("F20").Value = 91
$ws1.Range("F22").Value = 3256
$ws1.Range("F23").Value = 5598
$ws1.Range("F29").Value = 3201
$ws1.Range("F31").Value = 2402
$ws1.Range("F34").Value = 111
$ws1.Range("F36").Value = 245
$ws1.Range("F37").Value = 340
$ws1.Range("F38").Value = 97
$ws1.Range("F45").Value = 533

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 3769
$ws4.Range("F9").Value  = 112
$ws4.Range("F13").Value = 81
$ws4.Range("F16").Value = 895
$ws4.Range("G17").Value = 55
$ws4.Range("F18").Value = 226
$ws4.Range("F21").Value = 91
$ws4.Range("F23").Value = 3256
$ws4.Range("F24").Value = 5598
$ws4.Range("F30").Value = 3201
$ws4.Range("F32").Value = 2402
$ws4.Range("F37").Value = 245
$ws4.Range("F38").Value = 340
$ws4.Range("F39").Value = 97
$ws4.Range("F46").Value = 533

$wb.Save()
